$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo: "1) Djanko Unchained" -> "1) Django Unchained" ---
$ws.Range("I2").Value() = "1) Django Unchained"

# --- Remove the "More to add..." note from B12 and reset its row height ---
$ws.Range("B12").Clear()
$ws.Rows.Item(12).AutoFit()

# --- Row 13: Django Unchained Trailer ---
$rng = $ws.Range("A13")
$rng.Value() = "Django Unchained Trailer"
$ws.Hyperlinks.Add($rng, "https://www.youtube.com/watch?v=_xzQ7YVzoLk") | Out-Null
$rng.WrapText = $true
$ws.Range("C13").Value() = "YouTube - Movieclips Trailers"
$ws.Rows.Item(13).RowHeight = 30

# --- Row 14: Taxi Driver Trailer ---
$rng = $ws.Range("A14")
$rng.Value() = "Taxi Driver Trailer"
$ws.Hyperlinks.Add($rng, "https://www.youtube.com/watch?v=UeIBIz-0EOM") | Out-Null
$rng.WrapText = $true
$ws.Range("C14").Value() = "YouTube - Sony Pictures Home Entertainment"
$ws.Rows.Item(14).RowHeight = 45

# --- Row 15: Clockwork Orange Trailer ---
$rng = $ws.Range("A15")
$rng.Value() = "Clockwork Orange Trailer"
$ws.Hyperlinks.Add($rng, "https://www.youtube.com/watch?v=oMJBM1Et0ec") | Out-Null
$rng.WrapText = $true

# --- Row 16: American History X Trailer (typed before C15 to match shared-string order) ---
$rng = $ws.Range("A16")
$rng.Value() = "American History X Trailer"
$ws.Hyperlinks.Add($rng, "https://www.youtube.com/watch?v=qte9088cm3Q") | Out-Null
$rng.WrapText = $true

$ws.Range("C15").Value() = "YouTube - Warner Bros. Entertainment"
$ws.Range("C16").Value() = "YouTube - Movieclips Classic Trailers"
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 45

# --- Row 17: Full Metal Jacket Trailer ---
$rng = $ws.Range("A17")
$rng.Value() = "Full Metal Jacket Trailer"
$ws.Hyperlinks.Add($rng, "https://www.youtube.com/watch?v=g8JevAMv4-U") | Out-Null
$rng.WrapText = $true
$ws.Range("C17").Value() = "YouTube - SuperSubject20"
$ws.Rows.Item(17).RowHeight = 30

# --- Row 19-23: posters (row 18 intentionally left blank) ---
$rng = $ws.Range("A19")
$rng.Value() = "Rick and Morty Poster"
$ws.Hyperlinks.Add($rng, "https://live.staticflickr.com/65535/52145043728_rickandmorty.jpg") | Out-Null
$rng.WrapText = $true
$ws.Rows.Item(19).RowHeight = 30

$rng = $ws.Range("A20")
$rng.Value() = "Breaking Bad Poster"
$ws.Hyperlinks.Add($rng, "https://live.staticflickr.com/65535/52145043729_breakingbad.jpg") | Out-Null
$rng.WrapText = $true
$ws.Rows.Item(20).RowHeight = 30

$rng = $ws.Range("A21")
$rng.Value() = "Avatar: The Last Airbender Poster"
$ws.Hyperlinks.Add($rng, "https://live.staticflickr.com/65535/52145043730_avatartlab.jpg") | Out-Null
$rng.WrapText = $true
$ws.Rows.Item(21).RowHeight = 30

$rng = $ws.Range("A22")
$rng.Value() = "Lie To Me Poster"

$rng = $ws.Range("A23")
$rng.Value() = "Letterkenny Poster"
$ws.Hyperlinks.Add($rng, "https://live.staticflickr.com/65535/52145043732_letterkenny.jpg") | Out-Null
$rng.WrapText = $true

$rng = $ws.Range("A22")
$ws.Hyperlinks.Add($rng, "https://live.staticflickr.com/65535/52145043731_lietome.jpg") | Out-Null
$rng.WrapText = $true

# --- Update selection to reflect the saved view (scrolled down to the new rows) ---
$ws.Range("F20").Select()

Write-Host "Media list updated"
